$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1739130434782609
$ws.Range("C2").Value = 0.5852842809364549
$ws.Range("J2").Value = 0.02675585284280936
$ws.Range("P2").Value = 0.1237458193979933
$ws.Range("S2").Value = 0.0903010033444816
$ws.Range("B3").Value = 0.00546448087431694
$ws.Range("C3").Value = 0.02185792349726776
$ws.Range("J3").Value = 0.0273224043715847
$ws.Range("P3").Value = 0.7486338797814208
$ws.Range("S3").Value = 0.1967213114754098
$ws.Range("J4").Value = 0.07142857142857142
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2619047619047619
$ws.Range("B6").Value = 0.08225108225108226
$ws.Range("D6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.06060606060606061
$ws.Range("J6").Value = 0.2251082251082251
$ws.Range("O6").Value = 0.01298701298701299
$ws.Range("Q6").Value = 0.1774891774891775
$ws.Range("R6").Value = 0.08225108225108226
$ws.Range("S6").Value = 0.354978354978355
$ws.Range("B7").Value = 0.160427807486631
$ws.Range("D7").Value = 0.0213903743315508
$ws.Range("F7").Value = 0.05882352941176471
$ws.Range("J7").Value = 0.106951871657754
$ws.Range("O7").Value = 0.0106951871657754
$ws.Range("Q7").Value = 0.1443850267379679
$ws.Range("R7").Value = 0.09090909090909091
$ws.Range("S7").Value = 0.4064171122994653
$ws.Range("B8").Value = 0.09669211195928754
$ws.Range("D8").Value = 0.02035623409669211
$ws.Range("F8").Value = 0.06615776081424936
$ws.Range("J8").Value = 0.1272264631043257
$ws.Range("O8").Value = 0.0178117048346056
$ws.Range("Q8").Value = 0.1933842239185751
$ws.Range("R8").Value = 0.09923664122137404
$ws.Range("S8").Value = 0.3791348600508906
$ws.Range("B9").Value = 0.09852216748768473
$ws.Range("D9").Value = 0.01477832512315271
$ws.Range("E9").Value = 0.004926108374384237
$ws.Range("F9").Value = 0.06896551724137931
$ws.Range("J9").Value = 0.1133004926108374
$ws.Range("O9").Value = 0.009852216748768473
$ws.Range("Q9").Value = 0.1970443349753695
$ws.Range("R9").Value = 0.06403940886699508
$ws.Range("S9").Value = 0.4285714285714285
$ws.Range("B10").Value = 0.1062648691514671
$ws.Range("D10").Value = 0.02141157811260904
$ws.Range("F10").Value = 0.08643933386201427
$ws.Range("J10").Value = 0.1046788263283109
$ws.Range("O10").Value = 0.0126883425852498
$ws.Range("Q10").Value = 0.190325138778747
$ws.Range("R10").Value = 0.0816812053925456
$ws.Range("S10").Value = 0.3965107057890563
$ws.Range("G11").Value = 0.1638795986622074
$ws.Range("J11").Value = 0.08361204013377926
$ws.Range("K11").Value = 0.2307692307692308
$ws.Range("L11").Value = 0.5016722408026756
$ws.Range("S11").Value = 0.02006688963210702
$ws.Range("G12").Value = 0.7350993377483444
$ws.Range("J12").Value = 0.2119205298013245
$ws.Range("K12").Value = 0.01324503311258278
$ws.Range("L12").Value = 0.01324503311258278
$ws.Range("S12").Value = 0.02649006622516556
$ws.Range("G13").Value = 0.6086956521739131
$ws.Range("J13").Value = 0.3695652173913043
$ws.Range("S13").Value = 0.02173913043478261
$ws.Range("F15").Value = 0.02030456852791878
$ws.Range("H15").Value = 0.1421319796954315
$ws.Range("I15").Value = 0.09644670050761421
$ws.Range("J15").Value = 0.4060913705583756
$ws.Range("K15").Value = 0.04060913705583756
$ws.Range("M15").Value = 0.01522842639593909
$ws.Range("O15").Value = 0.04060913705583756
$ws.Range("S15").Value = 0.2385786802030457
$ws.Range("F16").Value = 0.005076142131979695
$ws.Range("H16").Value = 0.2030456852791878
$ws.Range("I16").Value = 0.09644670050761421
$ws.Range("J16").Value = 0.3959390862944163
$ws.Range("K16").Value = 0.08121827411167512
$ws.Range("M16").Value = 0.02030456852791878
$ws.Range("O16").Value = 0.04568527918781726
$ws.Range("S16").Value = 0.1522842639593909
$ws.Range("F17").Value = 0.01891252955082742
$ws.Range("H17").Value = 0.1607565011820331
$ws.Range("I17").Value = 0.1016548463356974
$ws.Range("J17").Value = 0.4444444444444444
$ws.Range("K17").Value = 0.08747044917257683
$ws.Range("M17").Value = 0.01654846335697399
$ws.Range("O17").Value = 0.0591016548463357
$ws.Range("S17").Value = 0.1111111111111111
$ws.Range("F18").Value = 0.01047120418848168
$ws.Range("H18").Value = 0.1727748691099476
$ws.Range("I18").Value = 0.1047120418848168
$ws.Range("J18").Value = 0.4712041884816754
$ws.Range("K18").Value = 0.1151832460732984
$ws.Range("M18").Value = 0.01570680628272251
$ws.Range("N18").Value = 0.005235602094240838
$ws.Range("O18").Value = 0.0418848167539267
$ws.Range("S18").Value = 0.06282722513089005
$ws.Range("F19").Value = 0.01295546558704453
$ws.Range("H19").Value = 0.1846153846153846
$ws.Range("I19").Value = 0.08421052631578947
$ws.Range("J19").Value = 0.3854251012145749
$ws.Range("K19").Value = 0.1165991902834008
$ws.Range("M19").Value = 0.02510121457489878
$ws.Range("N19").Value = 0.00242914979757085
$ws.Range("O19").Value = 0.07611336032388664
$ws.Range("S19").Value = 0.1125506072874494

Write-Host "Applied 109 cell updates"